# Ajout Buzzer + main
# - Renomme "Sheet2" en "UML"
# - Ajoute une colonne "Progres " a Table1 (feuille Objectifs)
# - Renseigne la progression pour la ligne "Detecteur de presence" et
#   "detecteur de porte ouverte/fermee" (passee a "En cours")
# - Ajoute "Porte" sur la feuille UML

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Objectifs")

# --- Table1: nouvelle colonne "Progres " --------------------------------
$tbl = $ws1.ListObjects.Item("Table1")
$tbl.ListColumns.Add() | Out-Null
$tbl.ListColumns.Item(7).Range.Cells.Item(1, 1).Value = "Progres "

# Valeurs de la nouvelle colonne
$ws1.Range("G2").Value = "Classe detection mouvement réalisée"

# --- Renommage de la 2e feuille + complement ----------------------------
$wsUML = $wb.Worksheets.Item("Sheet2")
$wsUML.Name = "UML"
$wsUML.Range("E2").Value = "Porte"

# --- Mise a jour de la ligne "detecteur de porte ouverte/fermee" --------
$ws1.Range("E3").Value = "En cours"
$ws1.Range("G3").Value = "à tester"

# --- Mise en forme de la nouvelle colonne -------------------------------
$ws1.Columns.Item(7).AutoFit() | Out-Null

# --- Selections actives (telles que laissees par l'auteur) -------------
$null = $wsUML.Range("E2").Select()
$null = $ws1.Range("G5").Select()
